{"js": "// Replace the text of the 15 ranked profit-list paragraphs according to the\n// new ranking: \"Stuffed Mushrooms\" becomes the new #1 (keeping the old #1's\n// value 7.35), every other dish shifts up one rank (keeping the value that\n// was already in that rank slot), \"Steak\" is re-inserted at rank 14 (value\n// 1.50, the old rank-15 value), and rank 15 (\"ggddg\") gets a new value 0.68.\n\nconst replacements = [\n  \"\ud83e\udd47 Stuffed Mushrooms profit: 7.35 \u20aa/min\",\n  \"\ud83e\udd48 Salad profit: 6.00 \u20aa/min\",\n  \"\ud83e\udd49 Brule Cream profit: 5.36 \u20aa/min\",\n  \"4) Pasta profit: 4.65 \u20aa/min\",\n  \"5) Pizza profit: 3.63 \u20aa/min\",\n  \"6) Krep profit: 3.56 \u20aa/min\",\n  \"7) Belgian Waffle profit: 3.20 \u20aa/min\",\n  \"8) Hamburger profit: 2.70 \u20aa/min\",\n  \"9) Empanadas profit: 2.68 \u20aa/min\",\n  \"10) Schnitzel profit: 2.48 \u20aa/min\",\n  \"11) Cake profit: 2.32 \u20aa/min\",\n  \"12) Roast profit: 2.00 \u20aa/min\",\n  \"13) Arancini profit: 1.81 \u20aa/min\",\n  \"14) Steak profit: 1.50 \u20aa/min\",\n  \"15) ggddg profit: 0.68 \u20aa/min\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the index of the first list paragraph (the one starting with the\n// gold-medal emoji) so we target the right run of paragraphs regardless of\n// what precedes them (title/subtitle).\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"\ud83e\udd47\") !== -1) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error(\"Could not locate the start of the profit list (\ud83e\udd47 paragraph).\");\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const para = paragraphs.items[startIndex + i];\n  para.getRange().insertText(replacements[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the text of the 15 ranked profit-list paragraphs according to the\n# new ranking: \"Stuffed Mushrooms\" becomes the new #1 (keeping the old #1's\n# value 7.35), every other dish shifts up one rank (keeping the value that\n# was already in that rank slot), \"Steak\" is re-inserted at rank 14 (value\n# 1.50, the old rank-15 value), and rank 15 (\"ggddg\") gets a new value 0.68.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    \"\ud83e\udd47 Stuffed Mushrooms profit: 7.35 \u20aa/min\",\n    \"\ud83e\udd48 Salad profit: 6.00 \u20aa/min\",\n    \"\ud83e\udd49 Brule Cream profit: 5.36 \u20aa/min\",\n    \"4) Pasta profit: 4.65 \u20aa/min\",\n    \"5) Pizza profit: 3.63 \u20aa/min\",\n    \"6) Krep profit: 3.56 \u20aa/min\",\n    \"7) Belgian Waffle profit: 3.20 \u20aa/min\",\n    \"8) Hamburger profit: 2.70 \u20aa/min\",\n    \"9) Empanadas profit: 2.68 \u20aa/min\",\n    \"10) Schnitzel profit: 2.48 \u20aa/min\",\n    \"11) Cake profit: 2.32 \u20aa/min\",\n    \"12) Roast profit: 2.00 \u20aa/min\",\n    \"13) Arancini profit: 1.81 \u20aa/min\",\n    \"14) Steak profit: 1.50 \u20aa/min\",\n    \"15) ggddg profit: 0.68 \u20aa/min\"\n)\n\n# Find the paragraph that starts the ranked list (the gold-medal emoji line)\n# so we target the right run of paragraphs regardless of what precedes them.\n$startIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*\ud83e\udd47*\") {\n        $startIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1) {\n    throw \"Could not locate the start of the profit list (\ud83e\udd47 paragraph).\"\n}\n\nfor ($j = 0; $j -lt $replacements.Count; $j++) {\n    $para = $d.Paragraphs($startIndex + $j)\n    $para.Range.Text = $replacements[$j]\n}\n"}
